$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 9492.727999999999
$ws.Range("I62").Value = 4807.3335
$ws.Range("J62").Value = 11249.75
$ws.Range("K62").Value = 4807.3335
$ws.Range("L62").Value = 11249.75
$ws.Range("M62").Value = -4183.3335
$ws.Range("N62").Value = -12497.75
$ws.Range("H65").Value = 9492.727999999999
$ws.Range("I65").Value = 4807.3335
$ws.Range("J65").Value = 11249.75
$ws.Range("K65").Value = 24036.6675
$ws.Range("L65").Value = 56248.75
$ws.Range("M65").Value = -20916.6675
$ws.Range("N65").Value = -62488.75
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("H87").Value = 94999.5
$ws.Range("J87").Value = 94999.5
$ws.Range("L87").Value = 94999.5
$ws.Range("N87").Value = -97495.5
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("H90").Value = 94999.5
$ws.Range("J90").Value = 94999.5
$ws.Range("L90").Value = 284998.5
$ws.Range("N90").Value = -297478.5
$ws.Range("H106").Value = 5380.7
$ws.Range("I106").Value = 5034.1113
$ws.Range("K106").Value = 5034.1113
$ws.Range("M106").Value = -4403.1113
$ws.Range("H111").Value = 1094
$ws.Range("J111").Value = 1500
$ws.Range("L111").Value = 4500
$ws.Range("N111").Value = -10634
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("M116").ClearContents()
$ws.Range("H118").Value = 916
$ws.Range("I118").Value = 282.33334
$ws.Range("J118").Value = 2056.6
$ws.Range("K118").Value = 847.0000200000001
$ws.Range("L118").Value = 6169.799999999999
$ws.Range("M118").Value = 809.9999799999999
$ws.Range("N118").Value = -9483.799999999999
$ws.Range("H125").Value = 2826.1
$ws.Range("I125").Value = 2479.5
$ws.Range("K125").Value = 22315.5
$ws.Range("M125").Value = -19855.5
$ws.Range("H138").Value = 2125.4614
$ws.Range("I138").Value = 1017.6667
$ws.Range("J138").Value = 3636.0908
$ws.Range("K138").Value = 3053.0001
$ws.Range("L138").Value = 10908.2724
$ws.Range("M138").Value = 2086.9999
$ws.Range("N138").Value = -21188.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4913.9
$ws.Range("I2").Value = 2962
$ws.Range("K2").Value = 2962
$ws.Range("M2").Value = -2849
$ws.Range("H74").Value = 8324.75
$ws.Range("I74").Value = 7250
$ws.Range("J74").Value = 9399.5
$ws.Range("K74").Value = 7250
$ws.Range("L74").Value = 9399.5
$ws.Range("M74").Value = -6376
$ws.Range("N74").Value = -11147.5
$ws.Range("H77").Value = 8324.75
$ws.Range("I77").Value = 7250
$ws.Range("J77").Value = 9399.5
$ws.Range("K77").Value = 36250
$ws.Range("L77").Value = 46997.5
$ws.Range("M77").Value = -31882
$ws.Range("N77").Value = -55733.5
$ws.Range("H116").Value = 4913.9
$ws.Range("I116").Value = 2962
$ws.Range("K116").Value = 2962
$ws.Range("M116").Value = -668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4913.9
$ws.Range("I3").Value = 2962
$ws.Range("K3").Value = 2962
$ws.Range("M3").Value = -2848
$ws.Range("H134").Value = 2723.818
$ws.Range("I134").Value = 2417.1052
$ws.Range("K134").Value = 7251.3156
$ws.Range("M134").Value = -4716.3156

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1099.8572
$ws.Range("I16").Value = 1099.8
$ws.Range("K16").Value = 1099.8
$ws.Range("M16").Value = -812.8
$ws.Range("H41").Value = 58096.668
$ws.Range("J41").Value = 58096.668
$ws.Range("L41").Value = 58096.668
$ws.Range("N41").Value = -58952.668
$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 97945
$ws.Range("J51").Value = 97945
$ws.Range("L51").Value = 97945
$ws.Range("N51").Value = -99417
$ws.Range("H59").Value = 34842
$ws.Range("J59").Value = 38947.25
$ws.Range("L59").Value = 38947.25
$ws.Range("N59").Value = -41237.25
$ws.Range("H61").Value = 97945
$ws.Range("J61").Value = 97945
$ws.Range("L61").Value = 97945
$ws.Range("N61").Value = -98641
$ws.Range("H113").Value = 1099.8572
$ws.Range("I113").Value = 1099.8
$ws.Range("K113").Value = 1099.8
$ws.Range("M113").Value = 1070.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 336.92856
$ws.Range("J23").Value = 384.45456
$ws.Range("L23").Value = 1153.36368
$ws.Range("N23").Value = -1623.36368
$ws.Range("H32").Value = 335.6
$ws.Range("I32").Value = 167.5
$ws.Range("K32").Value = 502.5
$ws.Range("M32").Value = -219.5
$ws.Range("H68").Value = 900
$ws.Range("J68").Value = 966.6667
$ws.Range("L68").Value = 2900.0001
$ws.Range("N68").Value = -4522.0001
$ws.Range("H71").Value = 900
$ws.Range("J71").Value = 966.6667
$ws.Range("L71").Value = 8700.0003
$ws.Range("N71").Value = -16812.0003
$ws.Range("H103").Value = 2174
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 2174
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 6522
$ws.Range("N103").Value = -8280
$ws.Range("M103").ClearContents()
$ws.Range("H114").Value = 2889.5
$ws.Range("J114").Value = 3987.3333
$ws.Range("L114").Value = 11961.9999
$ws.Range("N114").Value = -18469.9999
$ws.Range("H123").Value = 999
$ws.Range("I123").Value = 999
$ws.Range("K123").Value = 2997
$ws.Range("M123").Value = -547
$ws.Range("H129").Value = 1619.5
$ws.Range("I129").Value = 822.7143
$ws.Range("J129").Value = 3478.6667
$ws.Range("K129").Value = 2468.1429
$ws.Range("L129").Value = 10436.0001
$ws.Range("M129").Value = 2531.8571
$ws.Range("N129").Value = -20436.0001
$ws.Range("H131").Value = 1549
$ws.Range("I131").Value = 1417.3636
$ws.Range("J131").Value = 1790.3334
$ws.Range("K131").Value = 4252.0908
$ws.Range("L131").Value = 5371.0002
$ws.Range("M131").Value = 787.9092000000001
$ws.Range("N131").Value = -15451.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 4789.75
$ws.Range("I55").Value = 4954.8335
$ws.Range("K55").Value = 4954.8335
$ws.Range("M55").Value = -4627.8335
$ws.Range("H102").Value = 1124.7826
$ws.Range("I102").Value = 1124.7826
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1124.7826
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 497.2174
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 190
$ws.Range("I9").Value = 190
$ws.Range("K9").Value = 190
$ws.Range("M9").Value = 34
$ws.Range("H22").Value = 600
$ws.Range("J22").Value = 600
$ws.Range("L22").Value = 600
$ws.Range("N22").Value = -1190
$ws.Range("H27").Value = 600
$ws.Range("J27").Value = 600
$ws.Range("L27").Value = 600
$ws.Range("N27").Value = -814
$ws.Range("H46").Value = 6330.5
$ws.Range("I46").Value = 4271.1665
$ws.Range("J46").Value = 7875
$ws.Range("K46").Value = 4271.1665
$ws.Range("L46").Value = 7875
$ws.Range("M46").Value = -4083.1665
$ws.Range("N46").Value = -8251
$ws.Range("H68").Value = 8891.538
$ws.Range("I68").Value = 8000
$ws.Range("J68").Value = 9159
$ws.Range("K68").Value = 8000
$ws.Range("L68").Value = 9159
$ws.Range("M68").Value = -7251
$ws.Range("N68").Value = -10657
$ws.Range("H71").Value = 8891.538
$ws.Range("I71").Value = 8000
$ws.Range("J71").Value = 9159
$ws.Range("K71").Value = 40000
$ws.Range("L71").Value = 45795
$ws.Range("M71").Value = -36256
$ws.Range("N71").Value = -53283
$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -6450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7337.8
$ws.Range("I62").Value = 5599.5
$ws.Range("J62").Value = 8496.666999999999
$ws.Range("K62").Value = 5599.5
$ws.Range("L62").Value = 8496.666999999999
$ws.Range("M62").Value = -4975.5
$ws.Range("N62").Value = -9744.666999999999
$ws.Range("H65").Value = 7337.8
$ws.Range("I65").Value = 5599.5
$ws.Range("J65").Value = 8496.666999999999
$ws.Range("K65").Value = 27997.5
$ws.Range("L65").Value = 42483.335
$ws.Range("M65").Value = -24877.5
$ws.Range("N65").Value = -48723.335
